$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 109 (add SKU first, matching original authoring order) ---
$ws.Cells.Item(109, 2).Value = "8B-9WVK-ISVT"

# --- Row 107 FAMILIA (new shared string "Theramed") ---
$ws.Cells.Item(107, 6).Value = "Theramed"

# --- Row 106 TITULO ---
$ws.Cells.Item(106, 5).Value = "Sognare® Almohada Fussión 6 Pack Tamaño Estándar, Relleno Suave, Anti Ácaros y Lavable. Hipoalergénico, Anti Ácaros y Lavable. Cont. 6 Pzas."

# --- Row 106 ASIN ---
$ws.Cells.Item(106, 1).Value = "B0B7VF1S57"

# --- Row 107 ASIN ---
$ws.Cells.Item(107, 1).Value = "B0DFKVST8R"

# --- Row 107 TITULO ---
$ws.Cells.Item(107, 5).Value = "Sognare Theramed PRO Colchoneta de Masaje con 3 Terapias: Calor Infrarrojo, Vibración y Masaje Shiatsu. Alivia la Tensión, Rigidez y Dolor Corporal. 100 Noches de Garantía."

# --- Row 108 TITULO ---
$ws.Cells.Item(108, 5).Value = "Sognare Set 1 Cubre Colchón Individual + 1 Almohada Estandar Fussión Firme + 1 Almohada Suave. Hipoalergénico, Anti Ácaros y Lavable. Cont. 3 Pzas."

# --- Row 108 ASIN ---
$ws.Cells.Item(108, 1).Value = "B0B76D5HF2"

# --- Remaining FAMILIA DE PRODUCTO updates (reuse existing shared strings) ---
$ws.Cells.Item(106, 6).Value = "Almohada"
$ws.Cells.Item(108, 6).Value = "Cubre"

# --- Fill in the rest of new row 109 ---
$ws.Cells.Item(109, 1).Value = "B0B7QN1K82"
$ws.Cells.Item(109, 5).Value = "Sognare Set 1 Cubre Colchón King Size Extra Confort + 2 Almohadas Estándar Fussion Relleno Suave, 100% Algodón, Termorregulable, Hipoalergenico, Anti ácaros. Cont. 3 pzas."
$ws.Cells.Item(109, 6).Value = "Cubre"

# Match the row formatting used by the rest of the data rows
$ws.Rows.Item(109).RowHeight = 14.4

# --- Update the saved selection state to match the authored workbook ---
$ws.Range("B107").Select()
